$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 18-33: label "BD" -> "BDS" plus updated numeric values (B,C,D)
$rowsFull = @(
    @{ Row = 18;  Label = "BDS";                     B = 369; C = 1419; D = 5434 },
    @{ Row = 19;  Label = "BDS-PPO";                 B = 369; C = 1419; D = 5434 },
    @{ Row = 20;  Label = "BDS-TPF";                 B = 175; C = 696;  D = 2697 },
    @{ Row = 21;  Label = "BDS-PPO-TPF";             B = 175; C = 696;  D = 2697 },
    @{ Row = 22;  Label = "BDS-BRC";                 B = 53;  C = 123;  D = 306 },
    @{ Row = 23;  Label = "BDS-PPO-BRC";             B = 53;  C = 123;  D = 306 },
    @{ Row = 24;  Label = "BDS-BRC-TPF";             B = 65;  C = 165;  D = 528 },
    @{ Row = 25;  Label = "BDS-PPO-BRC-TPF";         B = 65;  C = 165;  D = 528 },
    @{ Row = 26;  Label = "BDS-GLM";                 B = 428; C = 1473; D = 6086 },
    @{ Row = 27;  Label = "BDS-PPO-GLM";             B = 428; C = 1473; D = 6086 },
    @{ Row = 28;  Label = "BDS-GLM-TPF";             B = 188; C = 741;  D = 2672 },
    @{ Row = 29;  Label = "BDS-PPO-GLM-TPF";         B = 188; C = 741;  D = 2672 },
    @{ Row = 30;  Label = "BDS-BRC-GLM";             B = 115; C = 295;  D = 1087 },
    @{ Row = 31;  Label = "BDS-PPO-BRC-GLM";         B = 115; C = 295;  D = 1087 },
    @{ Row = 32;  Label = "BDS-BRC-GLM-TPF";         B = 68;  C = 173;  D = 488 },
    @{ Row = 33;  Label = "BDS-PPO-BRC-GLM-TPF";     B = 68;  C = 173;  D = 488 }
)

foreach ($r in $rowsFull) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Label
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}

# Rows 50-65: only the label changes from "BD-..." to "BDS-..." (numeric values unchanged)
$labelOnly = @{
    50 = "BDS-JPS"
    51 = "BDS-PPO-JPS"
    52 = "BDS-TPF-JPS"
    53 = "BDS-PPO-TPF-JPS"
    54 = "BDS-BRC-JPS"
    55 = "BDS-PPO-BRC-JPS"
    56 = "BDS-BRC-TPF-JPS"
    57 = "BDS-PPO-BRC-TPF-JPS"
    58 = "BDS-GLM-JPS"
    59 = "BDS-PPO-GLM-JPS"
    60 = "BDS-GLM-TPF-JPS"
    61 = "BDS-PPO-GLM-TPF-JPS"
    62 = "BDS-BRC-GLM-JPS"
    63 = "BDS-PPO-BRC-GLM-JPS"
    64 = "BDS-BRC-GLM-TPF-JPS"
    65 = "BDS-PPO-BRC-GLM-TPF-JPS"
}

foreach ($row in $labelOnly.Keys) {
    $ws.Cells.Item($row, 1).Value = $labelOnly[$row]
}
